# Apply "Treinamento dos primeiros 150 tweets" edit:
# Fill column B (Relevancia label) for rows 2..151 on the "Treinamento" sheet
# and update the sheet view (zoom + selection) to match the author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

$values = @(1,0,1,1,1,1,0,1,1,0,0,1,1,0,0,0,1,1,1,0,1,0,0,1,1,1,0,1,1,0,1,0,1,1,1,1,0,0,0,0,1,0,1,1,1,0,1,1,0,0,1,0,1,0,1,1,0,1,1,1,1,1,1,1,1,1,0,1,0,1,0,0,1,1,1,0,0,0,1,1,1,1,0,0,1,0,0,1,1,1,1,0,0,1,1,1,0,1,0,0,1,1,1,1,0,1,1,1,1,1,1,0,0,1,1,1,1,0,0,0,1,0,1,1,0,0,0,1,1,1,0,1,1,0,1,1,0,1,1,0,1,1,0,0,1,1,1,1,0,0)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 2).Value = $v
    $row++
}

# Update the view: zoom to 70% and move the selection to B37
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("B37").Select()
